# "search details - done"
# Update property details on the "נכסים" (properties) sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("נכסים")

# Row 3 (Sapi's listing): rooms (H) and price (I) were placeholder text
# "test" - fill in the real search/listing numbers.
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 222

# Row 6 (Sunny's listing): update search count and property type; also
# normalize the other numeric fields on the row.
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 6
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 23
$ws.Range("K6").Value = "קרקע"
$ws.Range("O6").Value = 0
